$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value for every data row
# (rows 2 through 489). All of these cells currently contain 45203 and
# need to be updated to 45204.
$range = $ws.Range("C2:C489")
$range.Value = 45204
